# Rebuild Sheet1 with the new 4-column layout (index, codice, data, val)
# replacing the old 3-column layout (CODICE PUNTO, DATA, MISURA SOGGIACENZA [m]).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve the existing "header" style (bold font + thin border on all
# sides + centered/top alignment) by copying it out of A1 (which already
# carries that style) into a scratch cell before we wipe the sheet.
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Preserve the date-format style (custom numFmt "YYYY-MM-DD HH:MM:SS")
# by copying it out of B2 into a scratch cell too.
$ws.Range("B2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Wipe all existing content + formatting in the old used range.
$ws.Range("A1:C4").Clear()

# --- New headers: row 1 only has entries in B:D (A1 stays blank, matching
# the pandas-style index column with no header label).
$ws.Range("B1").Value = "codice"
$ws.Range("C1").Value = "data"
$ws.Range("D1").Value = "val"

# Apply the header style (bold/border/centered-top) to the new header cells.
$ws.Range("Z1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Data rows 2-4: A = positional index (0,1,2), B = codice, C = data
# (date serial, formatted), D = val (measurement).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "PO0170020U0002"
$ws.Range("C2").Value = 44350
$ws.Range("D2").Value = 73

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "PO0170020U0002"
$ws.Range("C3").Value = 44421
$ws.Range("D3").Value = 70.8

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "PO0170020U0002"
$ws.Range("C4").Value = 44475
$ws.Range("D4").Value = 69.72

# Apply the header-ish style (bold/border/centered-top) to the index column.
$ws.Range("Z1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Apply the date-format style to the "data" column.
$ws.Range("Z2").Copy()
$ws.Range("C2:C4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Clean up scratch cells used to transport styles.
$ws.Range("Z1:Z2").Clear()

# --- Misc sheet-level bits that changed in the diff.
$null = $ws.Range("A1").Select()
